# EPS v3.3.1 -> v3.4.2 update for "Trans Const Cost per Unit Cap Dist.xlsx"
# Source switches from the UT Austin Energy Institute (2016) report to the
# Americans for a Clean Energy Grid / Grid Strategies (2021) report, and the
# Data sheet's supporting calculation is rebuilt around the new source's
# MW-mile figures.

$wb = $excel.ActiveWorkbook

$about   = $wb.Worksheets.Item("About")
$data    = $wb.Worksheets.Item("Data")
$tccpucd = $wb.Worksheets.Item("TCCpUCD")

# ----- About sheet: new source citation & new CPI-adjustment ratio -----
$about.Range("B3").Value  = "Americans for a Clean Energy Grid and Grid Strategies"
$about.Range("B4").Value  = 2021
$about.Range("B5").Value  = "Transmission Projects Ready To Go: Plugging Into America's Untapped Renewable Resources"
$about.Range("B6").Value  = "https://cleanenergygrid.org/wp-content/uploads/2019/04/Transmission-Projects-Ready-to-Go-Final.pdf"
$about.Range("B7").Value  = "Pages 11-12"
$about.Range("A10").Value = "We adjust 2021 dollars to 2012 dollars using the following conversion factor:"
$about.Range("A11").Value = 0.84730412960844359

# ----- Data sheet: drop the old AVERAGE(...) helper row, rebuild with the -----
# ----- new source's $/MW-mile derivation (capacity-miles & total cost)   -----
$data.Range("A13:B13").ClearContents()

$data.Range("A1").Formula = "=17*10^6"
$data.Range("B1").Value   = "MW-miles"

$data.Range("A2").Formula = "=33*10^9"
$data.Range("B2").Value   = "USD"

$data.Range("A3").Formula = "=A2/A1"
$data.Range("B3").Value   = "$ / MW-mile"

$null = $data.Range("A6").Select()

# ----- TCCpUCD sheet: cost per unit capacity distance now references the -----
# ----- Data sheet's derived $/MW-mile value instead of the old average   -----
$tccpucd.Range("B2").Formula = "=Data!A3*About!A11"

# The old illustrative picture on the Data sheet is no longer part of the
# workbook.
if ($data.Shapes.Count -gt 0) {
    for ($i = $data.Shapes.Count; $i -ge 1; $i--) {
        $data.Shapes.Item($i).Delete()
    }
}

# Restore the "About" sheet as the active tab/selection (unchanged from
# before the edit) since selecting cells on other sheets above would
# otherwise shift the active tab.
$about.Activate()
$null = $about.Range("A11").Select()
